{"js": "// RBA 2.3 - Relat\u00f3rio e Email\n// Replace the placeholder tokens \"REW\"/\"QWREW\"/\"Rew\"/\"rew\" with\n// \"QWER\"/\"QWR\"/\"Qwer\"/\"qwer\" respectively, in both the document body\n// (the bold \"A QWREW,\" convocation target) and the page header\n// (diretoria / endere\u00e7o / CEP / Tel / Email lines).\n\n// --- 1. Body: bold \"QWREW\" -> \"QWR\" -------------------------------------\nconst body = context.document.body;\nconst bodyResults = body.search(\"QWREW\", { matchCase: true, matchWholeWord: true });\nbodyResults.load(\"items\");\nawait context.sync();\n\nfor (const result of bodyResults.items) {\n  result.insertText(\"QWR\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// --- 2. Header: REW -> QWER, QWREW -> QWR, Rew -> Qwer, rew -> qwer -----\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < sections.items.length; i++) {\n  const header = sections.items[i].getHeader(\"Primary\");\n\n  const replacements = [\n    [\"REW\", \"QWER\"],\n    [\"QWREW\", \"QWR\"],\n    [\"Rew\", \"Qwer\"],\n    [\"rew\", \"qwer\"],\n  ];\n\n  const searchResultsList = [];\n  for (const [find] of replacements) {\n    const results = header.search(find, { matchCase: true, matchWholeWord: true });\n    results.load(\"items\");\n    searchResultsList.push(results);\n  }\n  await context.sync();\n\n  for (let r = 0; r < replacements.length; r++) {\n    const [, replaceWith] = replacements[r];\n    const results = searchResultsList[r];\n    for (const item of results.items) {\n      item.insertText(replaceWith, Word.InsertLocation.replace);\n    }\n  }\n  await context.sync();\n}\n", "ps1": "# RBA 2.3 - Relat\u00f3rio e Email\n# Replace the placeholder tokens \"REW\"/\"QWREW\"/\"Rew\"/\"rew\" with\n# \"QWER\"/\"QWR\"/\"Qwer\"/\"qwer\" respectively, in both the document body\n# (the bold \"A QWREW,\" convocation target) and the page header\n# (diretoria / endere\u00e7o / CEP / Tel / Email lines).\n\n$d = $word.ActiveDocument\n\n# --- 1. Body: bold \"QWREW\" -> \"QWR\" --------------------------------------\n$bodyRange = $d.Content\n$bodyRange.Find.ClearFormatting()\n$bodyRange.Find.Text = \"QWREW\"\n$bodyRange.Find.MatchCase = $true\n$bodyRange.Find.MatchWholeWord = $true\n$bodyRange.Find.MatchWildcards = $false\n$bodyRange.Find.Forward = $true\n$bodyRange.Find.Wrap = 1\n$bodyRange.Find.Replacement.ClearFormatting()\n$bodyRange.Find.Replacement.Text = \"QWR\"\n[void]$bodyRange.Find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2)\n\n# --- 2. Header: REW -> QWER, QWREW -> QWR, Rew -> Qwer, rew -> qwer ------\nfor ($s = 1; $s -le $d.Sections.Count; $s++) {\n    $section = $d.Sections.Item($s)\n    $header = $section.Headers.Item(1)\n\n    $pairs = @(\n        @(\"REW\", \"QWER\"),\n        @(\"QWREW\", \"QWR\"),\n        @(\"Rew\", \"Qwer\"),\n        @(\"rew\", \"qwer\")\n    )\n\n    foreach ($pair in $pairs) {\n        $find = $pair[0]\n        $replace = $pair[1]\n\n        $hdrRange = $header.Range\n        $hdrRange.Find.ClearFormatting()\n        $hdrRange.Find.Text = $find\n        $hdrRange.Find.MatchCase = $true\n        $hdrRange.Find.MatchWholeWord = $true\n        $hdrRange.Find.MatchWildcards = $false\n        $hdrRange.Find.Forward = $true\n        $hdrRange.Find.Wrap = 1\n        $hdrRange.Find.Replacement.ClearFormatting()\n        $hdrRange.Find.Replacement.Text = $replace\n        [void]$hdrRange.Find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2)\n    }\n}\n"}
